$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.113852262496948
$ws.Range("B1").Value = 2.35284423828125
$ws.Range("C1").Value = 2.454925298690796
$ws.Range("D1").Value = 3.135828971862793
$ws.Range("E1").Value = 2.683386564254761
